$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Input")
$ws.Range("B2").Value = "HMP Fred"
$ws.Range("C2").Value = "Freds County Court"
$ws.Range("C2").Font.Name = "Arial"
$ws.Range("C2").Font.Size = 10
$ws.Range("B2").Select()
